# Apply updated KNN imputation results to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = -12.829
$ws.Range("E4").Value  = 13.195
$ws.Range("E5").Value  = 13.511
$ws.Range("C6").Value  = -12.445
$ws.Range("C7").Value  = -12.673
$ws.Range("E8").Value  = 13.468
$ws.Range("C16").Value = -11.638
$ws.Range("E16").Value = 13.164
$ws.Range("C20").Value = -13.041
$ws.Range("E22").Value = 13.302
